# Insert a new weekly price-report row for "Macroferia Regional de Talca -
# Repollo" right before the existing row 499. Excel shifts the former rows
# 499:561 down to 500:562 (dates/prices keep their original order), and the
# newly opened row 499 gets the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 499:561 down to 500:562, opening a blank row 499.
$ws.Rows("499:499").Insert()

# Fill the new row 499 with this week's record.
$ws.Cells.Item(499, 1).Value2  = 5
$ws.Cells.Item(499, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(499, 3).Value2  = "Maule"
$ws.Cells.Item(499, 4).Value2  = 45212
$ws.Cells.Item(499, 5).Value2  = 7
$ws.Cells.Item(499, 6).Value2  = 100112006
$ws.Cells.Item(499, 7).Value2  = "Repollo"
$ws.Cells.Item(499, 8).Value2  = "Crespo record"
$ws.Cells.Item(499, 9).Value2  = "Primera"
$ws.Cells.Item(499, 10).Value2 = 3000
$ws.Cells.Item(499, 11).Value2 = 800
$ws.Cells.Item(499, 12).Value2 = 800
$ws.Cells.Item(499, 13).Value2 = 800
$ws.Cells.Item(499, 14).Value2 = '$/unidad'
$ws.Cells.Item(499, 15).Value2 = "Región del Maule"
$ws.Cells.Item(499, 16).Value2 = 800
$ws.Cells.Item(499, 17).Value2 = 1
$ws.Cells.Item(499, 18).Value2 = "Hortaliza"
